$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: numpy tutorial post
$ws.Range("D6").Value = "[numpy tutorial] numpy에서 대각선 값 채우기 in python"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/numpy-tutorial-numpy%EC%97%90%EC%84%9C-%EB%8C%80%EA%B0%81%EC%84%A0-%EA%B0%92-%EC%B1%84%EC%9A%B0%EA%B8%B0-in-python"

# Row 12: tensorflow blog post
$ws.Range("D12").Value = "“혼자 공부하는 머신러닝+딥러닝”이 출간되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/01/19/%ed%98%bc%ec%9e%90-%ea%b3%b5%eb%b6%80%ed%95%98%eb%8a%94-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d%eb%94%a5%eb%9f%ac%eb%8b%9d%ec%9d%b4-%ec%b6%9c%ea%b0%84%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

# Row 40: insightCampus post
$newSpan = [char]0x3C + "span class=" + [char]0x22 + "kboard-default-new-notify" + [char]0x22 + [char]0x3E + "New" + [char]0x3C + "/span" + [char]0x3E
$tab = [char]0x9
$tabs16 = "$tab$tab$tab$tab$tab$tab$tab$tab$tab$tab$tab$tab$tab$tab$tab$tab"
$tabs8 = "$tab$tab$tab$tab$tab$tab$tab$tab"
$title = "2020년 머신러닝 프로젝트 Top 10"
$countSpan = [char]0x3C + "span class=" + [char]0x22 + "kboard-comments-count" + [char]0x22 + [char]0x3E + [char]0x3C + "/span" + [char]0x3E
$ws.Range("D40").Value = $newSpan + $tabs16 + $title + $tabs8 + $countSpan
$ws.Range("E40").Value = "https://insightcampus.co.kr/insightcommunity/?uid=12973&mod=document&pageid=1"
